# Replace the raw-Python "code paragraph" + two BodyText "output" paragraphs
# with a single syntax-highlighted SourceCode paragraph (code) followed by a
# single SourceCode paragraph holding the captured program output.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert a run of text at $pos (a collapsed-range offset) and return
# the new cursor position. Style is applied in a later pass (see below) so
# that <w:br/> runs end up on their own <w:r> with no rPr, matching how Word
# itself serializes line breaks that were never explicitly (re)styled.
#
# $segments must be an ArrayList (not a plain @() literal) of 3-element
# arrays @(kind, text, styleId) -- using ArrayList.Add() avoids the
# nested-array flattening surprises plain array literals run into once they
# cross a function-parameter boundary.
# ---------------------------------------------------------------------------
function Insert-Segments($startPos, $segments) {
    $pos = $startPos
    $ranges = New-Object System.Collections.ArrayList
    foreach ($seg in $segments) {
        $kind = $seg[0]
        $text = $seg[1]
        $style = $seg[2]
        $r = $d.Range($pos, $pos)
        if ($kind -eq "BREAK") {
            $r.InsertAfter([char]11)
        } else {
            $r.InsertAfter($text)
        }
        [void]$ranges.Add(@($pos, $r.End, $style))
        $pos = $r.End
    }
    # Apply character styles back-to-front so an un-styled <w:br/> segment
    # never inherits the "current" formatting left behind by a later
    # style assignment.
    for ($i = $ranges.Count - 1; $i -ge 0; $i--) {
        $seg = $ranges[$i]
        if ($seg[2]) {
            $rr = $d.Range($seg[0], $seg[1])
            $rr.Style = $seg[2]
        }
    }
    return $pos
}

# ---------------------------------------------------------------------------
# Segment lists describing the final two paragraphs, derived run-by-run from
# the target markup: (kind, text, characterStyleId)
# ---------------------------------------------------------------------------
$segmentsCode = New-Object System.Collections.ArrayList
[void]$segmentsCode.Add(@("TEXT", "import", "ImportTok"))
[void]$segmentsCode.Add(@("TEXT", " platform", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "import", "ImportTok"))
[void]$segmentsCode.Add(@("TEXT", " psutil", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "print", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "`"Operating System:`"", "StringTok"))
[void]$segmentsCode.Add(@("TEXT", ", platform.system(), platform.release())", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "print", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "`"Python Version:`"", "StringTok"))
[void]$segmentsCode.Add(@("TEXT", ", platform.python_version())", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "print", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "`"Machine:`"", "StringTok"))
[void]$segmentsCode.Add(@("TEXT", ", platform.machine())", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "print", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "`"Processor:`"", "StringTok"))
[void]$segmentsCode.Add(@("TEXT", ", platform.processor())", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "mem ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "=", "OperatorTok"))
[void]$segmentsCode.Add(@("TEXT", " psutil.virtual_memory()", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "print", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "`"Total Memory (MB):`"", "StringTok"))
[void]$segmentsCode.Add(@("TEXT", ", ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "round", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(mem.total ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "/", "OperatorTok"))
[void]$segmentsCode.Add(@("TEXT", " (", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "1024", "DecValTok"))
[void]$segmentsCode.Add(@("TEXT", "**", "OperatorTok"))
[void]$segmentsCode.Add(@("TEXT", "2", "DecValTok"))
[void]$segmentsCode.Add(@("TEXT", "), ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "2", "DecValTok"))
[void]$segmentsCode.Add(@("TEXT", "))", "NormalTok"))
[void]$segmentsCode.Add(@("BREAK", $null, $null))
[void]$segmentsCode.Add(@("TEXT", "print", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "`"Available Memory (MB):`"", "StringTok"))
[void]$segmentsCode.Add(@("TEXT", ", ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "round", "BuiltInTok"))
[void]$segmentsCode.Add(@("TEXT", "(mem.available ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "/", "OperatorTok"))
[void]$segmentsCode.Add(@("TEXT", " (", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "1024", "DecValTok"))
[void]$segmentsCode.Add(@("TEXT", "**", "OperatorTok"))
[void]$segmentsCode.Add(@("TEXT", "2", "DecValTok"))
[void]$segmentsCode.Add(@("TEXT", "), ", "NormalTok"))
[void]$segmentsCode.Add(@("TEXT", "2", "DecValTok"))
[void]$segmentsCode.Add(@("TEXT", "))", "NormalTok"))

$segmentsOutput = New-Object System.Collections.ArrayList
[void]$segmentsOutput.Add(@("TEXT", "Operating System: Linux 6.14.0-1011-aws", "VerbatimChar"))
[void]$segmentsOutput.Add(@("BREAK", $null, $null))
[void]$segmentsOutput.Add(@("TEXT", "Python Version: 3.12.3", "VerbatimChar"))
[void]$segmentsOutput.Add(@("BREAK", $null, $null))
[void]$segmentsOutput.Add(@("TEXT", "Machine: x86_64", "VerbatimChar"))
[void]$segmentsOutput.Add(@("BREAK", $null, $null))
[void]$segmentsOutput.Add(@("TEXT", "Processor: x86_64", "VerbatimChar"))
[void]$segmentsOutput.Add(@("BREAK", $null, $null))
[void]$segmentsOutput.Add(@("TEXT", "Total Memory (MB): 7820.98", "VerbatimChar"))
[void]$segmentsOutput.Add(@("BREAK", $null, $null))
[void]$segmentsOutput.Add(@("TEXT", "Available Memory (MB): 6027.11", "VerbatimChar"))

# ---------------------------------------------------------------------------
# 1. Locate the three source paragraphs: "First Paragraph" (import lines)
#    followed by two "Body Text" paragraphs (the print()/mem lines).
# ---------------------------------------------------------------------------
$codeParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Style.NameLocal -eq "First Paragraph") {
        $codeParaIndex = $i
        break
    }
}

$pCode = $d.Paragraphs.Item($codeParaIndex)
$pPrint = $d.Paragraphs.Item($codeParaIndex + 1)
$pMem = $d.Paragraphs.Item($codeParaIndex + 2)

# ---------------------------------------------------------------------------
# 2. Delete the two trailing Body Text paragraphs (their text is being
#    absorbed into the rewritten code paragraph / replaced by the output
#    paragraph).
# ---------------------------------------------------------------------------
$deleteRange = $d.Range($pPrint.Range.Start, $pMem.Range.End)
$deleteRange.Delete()

# ---------------------------------------------------------------------------
# 3. Rebuild the code paragraph from $segmentsCode. The new runs are
#    inserted *before* the old "import platform import psutil" text (so the
#    paragraph is never emptied out, which would otherwise collapse/merge
#    it with its neighbour), and the stale original text is deleted
#    afterwards from its shifted position. Finally the paragraph itself is
#    restyled to "Source Code".
# ---------------------------------------------------------------------------
$pCode = $d.Paragraphs.Item($codeParaIndex)
$oldStart = $pCode.Range.Start
$oldEnd = $pCode.Range.End
$oldLen = $oldEnd - $oldStart

$newEnd = Insert-Segments $oldStart $segmentsCode
$newLen = $newEnd - $oldStart

$staleRange = $d.Range($oldStart + $newLen, $oldStart + $newLen + $oldLen)
$staleRange.Delete()

$pCode = $d.Paragraphs.Item($codeParaIndex)
$pCode.Range.Style = "SourceCode"

# ---------------------------------------------------------------------------
# 4. Insert a new paragraph after the code paragraph to hold the captured
#    program output, style it "Source Code", and populate it from
#    $segmentsOutput.
# ---------------------------------------------------------------------------
$pCode = $d.Paragraphs.Item($codeParaIndex)
$pCode.Range.InsertParagraphAfter() | Out-Null

$pOutput = $d.Paragraphs.Item($codeParaIndex + 1)
$pOutput.Range.Style = "SourceCode"

Insert-Segments $pOutput.Range.Start $segmentsOutput | Out-Null
